$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> new D (Price) / E (Volume 1h) values, exactly as they must be stored
# (both columns hold plain text in the source data, even when the text looks
# like a number).
$updates = @{
    2  = @{ D = "63.951.04"; E = "  +0.08%  " }
    3  = @{ D = "3.141.23";  E = "  +0.62%  " }
    4  = @{ E = "  +0.10%  " }
    5  = @{ D = "590.38";    E = "  +0.77%  " }
    6  = @{ D = "145.59";    E = "  -0.56%  " }
    7  = @{ E = "  +0.05%  " }
    8  = @{ D = "3.135.26";  E = "  +0.61%  " }
    9  = @{ E = "  -0.45%  " }
    10 = @{ E = "  -0.82%  " }
    11 = @{ D = "5.91";      E = "  +2.89%  " }
    12 = @{ E = "  -1.74%  " }
    13 = @{ E = "  -2.76%  " }
    14 = @{ D = "37.35";     E = "  +0.66%  " }
    15 = @{ D = "3.663.08";  E = "  +0.60%  " }
    16 = @{ E = "  -1.33%  " }
    17 = @{ D = "7.36";      E = "  +2.84%  " }
    18 = @{ D = "63.812.88"; E = "  +0.00%  " }
    19 = @{ D = "3.138.61";  E = "  +0.44%  " }
    20 = @{ D = "468.53";    E = "  +0.64%  " }
    21 = @{ E = "  +0.35%  " }
    22 = @{ E = "  +0.12%  " }
    23 = @{ E = "  -0.15%  " }
    24 = @{ D = "12.99";     E = "  -1.34%  " }
    25 = @{ D = "81.62";     E = "  -0.59%  " }
    26 = @{ D = "2.31";      E = "  +6.69%  " }
    27 = @{ E = "  +0.10%  " }
    28 = @{ D = "9.81";      E = "  +9.71%  " }
    29 = @{ D = "7.43";      E = "  +8.34%  " }
    30 = @{ E = "  +0.15%  " }
    31 = @{ E = "  +0.29%  " }
    32 = @{ E = "  +0.10%  " }
    33 = @{ D = "27.73";     E = "  +2.55%  " }
    34 = @{ E = "  +0.60%  " }
    35 = @{ D = "0.0₃0844";  E = "  -3.90%  " }
    36 = @{ D = "1.07";      E = "  +1.19%  " }
    37 = @{ E = "  +1.28%  " }
    38 = @{ E = "  -3.03%  " }
    39 = @{ D = "3.20";      E = "  -6.35%  " }
    40 = @{ E = "  +1.09%  " }
    41 = @{ D = "9.37";      E = "  +7.81%  " }
    42 = @{ D = "453.84";    E = "  +0.49%  " }
    43 = @{ D = "0.291";     E = "  +4.83%  " }
    44 = @{ E = "  -0.13%  " }
    45 = @{ D = "2.916.91";  E = "  +1.02%  " }
    46 = @{ D = "39.63";     E = "  +10.44%  " }
    47 = @{ E = "  -2.90%  " }
    48 = @{ D = "133.12";    E = "  +6.79%  " }
    50 = @{ E = "  +2.74%  " }
    51 = @{ E = "  -0.59%  " }
}

foreach ($row in $updates.Keys) {
    $rowData = $updates[$row]

    if ($rowData.ContainsKey("D")) {
        $dCell = $ws.Range("D$row")
        $newVal = $rowData["D"]

        # Decimal-looking strings (e.g. "590.38") get auto-coerced to a
        # number by Excel's normal type inference. Force text storage via a
        # Text number format, write the value, then restore the cell to its
        # original (default/general) style so no stray formatting lingers.
        $isNumericLooking = $newVal -match '^[+-]?[0-9]+(\.[0-9]+)?$'
        if ($isNumericLooking) {
            $dCell.NumberFormat = "@"
            $dCell.Value = $newVal
            $dCell.NumberFormat = "General"
            $dCell.Style = "Normal"
        } else {
            $dCell.Value = $newVal
        }
    }

    if ($rowData.ContainsKey("E")) {
        $ws.Range("E$row").Value = $rowData["E"]
    }
}
